# "Colocando header nos graficos"
# Adds a header label in column A (row 1) of each data sheet, strips the
# bold/border style from the row-label cells below it, fixes a handful of
# accented-Portuguese typos, drops the "Teto" row from the emissions sheet,
# and refreshes the header/values on the cost sheet.

$wb = $excel.ActiveWorkbook

function Set-HeaderCell {
    param($ws, [string]$text)
    # Give A1 the same look (bold / border / centered) as the rest of row 1
    # by copying the formatting from the neighboring B1 header cell, then
    # writing the label text into it.
    $ws.Range("B1").Copy() | Out-Null
    $ws.Range("A1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range("A1").Value = $text
}

function Clear-LabelStyle {
    param($ws, [string]$addr)
    # Row labels below the header used to share the bold/border style (s=1);
    # now only the header row keeps it.
    $ws.Range($addr).Style = "Normal"
}

# ---------------------------------------------------------------------
# Sheets 1-4: "Potencia Acumulada", "Geracao Periodo Medio",
# "Atendimento a Ponta", "Potencia Incremental" all share the same layout.
# ---------------------------------------------------------------------
for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)

    Set-HeaderCell $ws "Fonte/Tecnologia"

    Clear-LabelStyle $ws "A2"
    $ws.Range("A2").Value = "Hidro"

    Clear-LabelStyle $ws "A3"
    $ws.Range("A3").Value = "Gás Natural"

    Clear-LabelStyle $ws "A4"
    $ws.Range("A4").Value = "Carvão"

    Clear-LabelStyle $ws "A5"
    $ws.Range("A5").Value = "Nuclear"

    Clear-LabelStyle $ws "A6"
    $ws.Range("A6").Value = "Óleos Comb"

    Clear-LabelStyle $ws "A7"
    $ws.Range("A7").Value = "Biomassa"

    Clear-LabelStyle $ws "A8"
    $ws.Range("A8").Value = "Eólica"

    Clear-LabelStyle $ws "A9"
    $ws.Range("A9").Value = "Solar"

    Clear-LabelStyle $ws "A10"
    $ws.Range("A10").Value = "Outros"

    Clear-LabelStyle $ws "A11"
    $ws.Range("A11").Value = "Pot. Compl."

    Clear-LabelStyle $ws "A12"
    $ws.Range("A12").Value = "GD"
}

# ---------------------------------------------------------------------
# Sheet 5: "Emissoes Totais (MtCO2eq)"
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

Set-HeaderCell $ws5 "Período"

Clear-LabelStyle $ws5 "A2"
$ws5.Range("A2").Value = "P.Médio"

Clear-LabelStyle $ws5 "A3"
$ws5.Range("A3").Value = "P.Crítico"

# Row 4 ("Teto") is removed entirely.
$ws5.Rows.Item(4).Delete()

# ---------------------------------------------------------------------
# Sheet 6: "Custo Total (bilhoes de R$)"
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

Set-HeaderCell $ws6 "Tipo Expansão"

# B1 switches from the label "Custo" to the year "2015" (same s=1 style).
$ws6.Range("B1").Value = 2015

Clear-LabelStyle $ws6 "A2"
$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("B2").Value = 588

Clear-LabelStyle $ws6 "A3"
$ws6.Range("A3").Value = "Expansão por GD"
$ws6.Range("B3").Value = 99

Write-Host "Headers added and labels normalized."
